$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9263527393341064
$ws.Range("B1").Value = 1.557147979736328
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.920924782752991
$ws.Range("E1").Value = 1.271753668785095
